$wb = $excel.ActiveWorkbook

# ===== Sheet: TRX =====
$ws = $wb.Worksheets.Item("TRX")

# New column header AK1 (date 20220218), matching style of existing header cells
$ws.Cells.Item(1, 37).NumberFormat = "@"
$ws.Cells.Item(1, 37).Font.Bold = $true
$ws.Cells.Item(1, 37).HorizontalAlignment = -4108
$ws.Cells.Item(1, 37).VerticalAlignment = -4160
$ws.Cells.Item(1, 37).Borders.LineStyle = 1
$ws.Cells.Item(1, 37).Value = "20220218"

# New row 14: period 20210903-20210910
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Font.Bold = $true
$ws.Cells.Item(14, 1).HorizontalAlignment = -4108
$ws.Cells.Item(14, 1).VerticalAlignment = -4160
$ws.Cells.Item(14, 1).Borders.LineStyle = 1
$ws.Cells.Item(14, 1).Value = "20210903-20210910"
$ws.Cells.Item(14, 14).Value = 422993.89125
$ws.Cells.Item(14, 15).Value = 422993.89125
$ws.Cells.Item(14, 16).Value = 422993.89125
$ws.Cells.Item(14, 17).Value = 422993.89125
$ws.Cells.Item(14, 18).Value = 422993.89125
$ws.Cells.Item(14, 19).Value = 422993.89125
$ws.Cells.Item(14, 20).Value = 422993.89125
$ws.Cells.Item(14, 21).Value = 422993.89125
$ws.Cells.Item(14, 22).Value = 422993.89125
$ws.Cells.Item(14, 23).Value = 422993.89125
$ws.Cells.Item(14, 24).Value = 422993.89125
$ws.Cells.Item(14, 25).Value = 422993.89125
$ws.Cells.Item(14, 26).Value = 422993.89125
$ws.Cells.Item(14, 27).Value = 422993.89125
$ws.Cells.Item(14, 28).Value = 422993.89125
$ws.Cells.Item(14, 29).Value = 422993.89125
$ws.Cells.Item(14, 30).Value = 422993.89125
$ws.Cells.Item(14, 31).Value = 422993.89125
$ws.Cells.Item(14, 32).Value = 422993.89125
$ws.Cells.Item(14, 33).Value = 422993.89125
$ws.Cells.Item(14, 34).Value = 422993.89125
$ws.Cells.Item(14, 35).Value = 422993.89125
$ws.Cells.Item(14, 36).Value = 422993.89125
$ws.Cells.Item(14, 37).Value = 422993.89125

# Row 15 (previously row 14): SUM row, shifted down, with new cumulative totals
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Font.Bold = $true
$ws.Cells.Item(15, 1).HorizontalAlignment = -4108
$ws.Cells.Item(15, 1).VerticalAlignment = -4160
$ws.Cells.Item(15, 1).Borders.LineStyle = 1
$ws.Cells.Item(15, 1).Value = "SUM"
$ws.Cells.Item(15, 2).Value = 939796.0891666667
$ws.Cells.Item(15, 3).Value = 1879592.178333333
$ws.Cells.Item(15, 4).Value = 2819388.2675
$ws.Cells.Item(15, 5).Value = 3759184.356666667
$ws.Cells.Item(15, 6).Value = 4801479.284583334
$ws.Cells.Item(15, 7).Value = 5843774.2125
$ws.Cells.Item(15, 8).Value = 6886069.140416667
$ws.Cells.Item(15, 9).Value = 7928364.068333333
$ws.Cells.Item(15, 10).Value = 8765438.304583333
$ws.Cells.Item(15, 11).Value = 9602512.540833334
$ws.Cells.Item(15, 12).Value = 10439586.77708333
$ws.Cells.Item(15, 13).Value = 11276661.01333333
$ws.Cells.Item(15, 14).Value = 11699654.904583333
$ws.Cells.Item(15, 15).Value = 11699654.904583333
$ws.Cells.Item(15, 16).Value = 11699654.904583333
$ws.Cells.Item(15, 17).Value = 11699654.904583333
$ws.Cells.Item(15, 18).Value = 11699654.904583333
$ws.Cells.Item(15, 19).Value = 11699654.904583333
$ws.Cells.Item(15, 20).Value = 11699654.904583333
$ws.Cells.Item(15, 21).Value = 11699654.904583333
$ws.Cells.Item(15, 22).Value = 11699654.904583333
$ws.Cells.Item(15, 23).Value = 11699654.904583333
$ws.Cells.Item(15, 24).Value = 11699654.904583333
$ws.Cells.Item(15, 25).Value = 11699654.904583333
$ws.Cells.Item(15, 26).Value = 10759858.815416668
$ws.Cells.Item(15, 27).Value = 9820062.72625
$ws.Cells.Item(15, 28).Value = 8880266.637083335
$ws.Cells.Item(15, 29).Value = 7940470.5479166685
$ws.Cells.Item(15, 30).Value = 6898175.620000002
$ws.Cells.Item(15, 31).Value = 5855880.692083335
$ws.Cells.Item(15, 32).Value = 4813585.764166667
$ws.Cells.Item(15, 33).Value = 3771290.83625
$ws.Cells.Item(15, 34).Value = 2934216.6
$ws.Cells.Item(15, 35).Value = 2097142.36375
$ws.Cells.Item(15, 36).Value = 1260068.1275
$ws.Cells.Item(15, 37).Value = 422993.89125

# ===== Sheet: JST =====
$ws = $wb.Worksheets.Item("JST")

# New column header AK1 (date 20220218), matching style of existing header cells
$ws.Cells.Item(1, 37).NumberFormat = "@"
$ws.Cells.Item(1, 37).Font.Bold = $true
$ws.Cells.Item(1, 37).HorizontalAlignment = -4108
$ws.Cells.Item(1, 37).VerticalAlignment = -4160
$ws.Cells.Item(1, 37).Borders.LineStyle = 1
$ws.Cells.Item(1, 37).Value = "20220218"

# New row 14: period 20210903-20210910
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Font.Bold = $true
$ws.Cells.Item(14, 1).HorizontalAlignment = -4108
$ws.Cells.Item(14, 1).VerticalAlignment = -4160
$ws.Cells.Item(14, 1).Borders.LineStyle = 1
$ws.Cells.Item(14, 1).Value = "20210903-20210910"
$ws.Cells.Item(14, 14).Value = 701978.0920833334
$ws.Cells.Item(14, 15).Value = 701978.0920833334
$ws.Cells.Item(14, 16).Value = 701978.0920833334
$ws.Cells.Item(14, 17).Value = 701978.0920833334
$ws.Cells.Item(14, 18).Value = 701978.0920833334
$ws.Cells.Item(14, 19).Value = 701978.0920833334
$ws.Cells.Item(14, 20).Value = 701978.0920833334
$ws.Cells.Item(14, 21).Value = 701978.0920833334
$ws.Cells.Item(14, 22).Value = 701978.0920833334
$ws.Cells.Item(14, 23).Value = 701978.0920833334
$ws.Cells.Item(14, 24).Value = 701978.0920833334
$ws.Cells.Item(14, 25).Value = 701978.0920833334
$ws.Cells.Item(14, 26).Value = 701978.0920833334
$ws.Cells.Item(14, 27).Value = 701978.0920833334
$ws.Cells.Item(14, 28).Value = 701978.0920833334
$ws.Cells.Item(14, 29).Value = 701978.0920833334
$ws.Cells.Item(14, 30).Value = 701978.0920833334
$ws.Cells.Item(14, 31).Value = 701978.0920833334
$ws.Cells.Item(14, 32).Value = 701978.0920833334
$ws.Cells.Item(14, 33).Value = 701978.0920833334
$ws.Cells.Item(14, 34).Value = 701978.0920833334
$ws.Cells.Item(14, 35).Value = 701978.0920833334
$ws.Cells.Item(14, 36).Value = 701978.0920833334
$ws.Cells.Item(14, 37).Value = 701978.0920833334

# Row 15 (previously row 14): SUM row, shifted down, with new cumulative totals
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Font.Bold = $true
$ws.Cells.Item(15, 1).HorizontalAlignment = -4108
$ws.Cells.Item(15, 1).VerticalAlignment = -4160
$ws.Cells.Item(15, 1).Borders.LineStyle = 1
$ws.Cells.Item(15, 1).Value = "SUM"
$ws.Cells.Item(15, 2).Value = 942083.3333333334
$ws.Cells.Item(15, 3).Value = 1829625
$ws.Cells.Item(15, 4).Value = 2717166.666666667
$ws.Cells.Item(15, 5).Value = 3604708.333333333
$ws.Cells.Item(15, 6).Value = 4477375
$ws.Cells.Item(15, 7).Value = 5350041.666666667
$ws.Cells.Item(15, 8).Value = 6222708.333333334
$ws.Cells.Item(15, 9).Value = 7095375.000000001
$ws.Cells.Item(15, 10).Value = 8022583.333333334
$ws.Cells.Item(15, 11).Value = 9004333.333333334
$ws.Cells.Item(15, 12).Value = 9986083.333333334
$ws.Cells.Item(15, 13).Value = 10967833.33333333
$ws.Cells.Item(15, 14).Value = 11669811.425416667
$ws.Cells.Item(15, 15).Value = 11669811.425416667
$ws.Cells.Item(15, 16).Value = 11669811.425416667
$ws.Cells.Item(15, 17).Value = 11669811.425416667
$ws.Cells.Item(15, 18).Value = 11669811.425416667
$ws.Cells.Item(15, 19).Value = 11669811.425416667
$ws.Cells.Item(15, 20).Value = 11669811.425416667
$ws.Cells.Item(15, 21).Value = 11669811.425416667
$ws.Cells.Item(15, 22).Value = 11669811.425416667
$ws.Cells.Item(15, 23).Value = 11669811.425416667
$ws.Cells.Item(15, 24).Value = 11669811.425416667
$ws.Cells.Item(15, 25).Value = 11669811.425416667
$ws.Cells.Item(15, 26).Value = 10727728.092083333
$ws.Cells.Item(15, 27).Value = 9840186.425416665
$ws.Cells.Item(15, 28).Value = 8952644.75875
$ws.Cells.Item(15, 29).Value = 8065103.092083333
$ws.Cells.Item(15, 30).Value = 7192436.425416667
$ws.Cells.Item(15, 31).Value = 6319769.758749999
$ws.Cells.Item(15, 32).Value = 5447103.092083333
$ws.Cells.Item(15, 33).Value = 4574436.425416667
$ws.Cells.Item(15, 34).Value = 3647228.0920833335
$ws.Cells.Item(15, 35).Value = 2665478.0920833335
$ws.Cells.Item(15, 36).Value = 1683728.0920833335
$ws.Cells.Item(15, 37).Value = 701978.0920833334

# ===== Sheet: WBTT =====
$ws = $wb.Worksheets.Item("WBTT")

# New column header AK1 (date 20220218), matching style of existing header cells
$ws.Cells.Item(1, 37).NumberFormat = "@"
$ws.Cells.Item(1, 37).Font.Bold = $true
$ws.Cells.Item(1, 37).HorizontalAlignment = -4108
$ws.Cells.Item(1, 37).VerticalAlignment = -4160
$ws.Cells.Item(1, 37).Borders.LineStyle = 1
$ws.Cells.Item(1, 37).Value = "20220218"

# New row 14: period 20210903-20210910
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Font.Bold = $true
$ws.Cells.Item(14, 1).HorizontalAlignment = -4108
$ws.Cells.Item(14, 1).VerticalAlignment = -4160
$ws.Cells.Item(14, 1).Borders.LineStyle = 1
$ws.Cells.Item(14, 1).Value = "20210903-20210910"
$ws.Cells.Item(14, 14).Value = 1513058.30375
$ws.Cells.Item(14, 15).Value = 1513058.30375
$ws.Cells.Item(14, 16).Value = 1513058.30375
$ws.Cells.Item(14, 17).Value = 1513058.30375
$ws.Cells.Item(14, 18).Value = 1513058.30375
$ws.Cells.Item(14, 19).Value = 1513058.30375
$ws.Cells.Item(14, 20).Value = 1513058.30375
$ws.Cells.Item(14, 21).Value = 1513058.30375
$ws.Cells.Item(14, 22).Value = 1513058.30375
$ws.Cells.Item(14, 23).Value = 1513058.30375
$ws.Cells.Item(14, 24).Value = 1513058.30375
$ws.Cells.Item(14, 25).Value = 1513058.30375
$ws.Cells.Item(14, 26).Value = 1513058.30375
$ws.Cells.Item(14, 27).Value = 1513058.30375
$ws.Cells.Item(14, 28).Value = 1513058.30375
$ws.Cells.Item(14, 29).Value = 1513058.30375
$ws.Cells.Item(14, 30).Value = 1513058.30375
$ws.Cells.Item(14, 31).Value = 1513058.30375
$ws.Cells.Item(14, 32).Value = 1513058.30375
$ws.Cells.Item(14, 33).Value = 1513058.30375
$ws.Cells.Item(14, 34).Value = 1513058.30375
$ws.Cells.Item(14, 35).Value = 1513058.30375
$ws.Cells.Item(14, 36).Value = 1513058.30375
$ws.Cells.Item(14, 37).Value = 1513058.30375

# Row 15 (previously row 14): SUM row, shifted down, with new cumulative totals
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Font.Bold = $true
$ws.Cells.Item(15, 1).HorizontalAlignment = -4108
$ws.Cells.Item(15, 1).VerticalAlignment = -4160
$ws.Cells.Item(15, 1).Borders.LineStyle = 1
$ws.Cells.Item(15, 1).Value = "SUM"
$ws.Cells.Item(15, 2).Value = 2022976.932083334
$ws.Cells.Item(15, 3).Value = 4045953.864166667
$ws.Cells.Item(15, 4).Value = 6068930.796250001
$ws.Cells.Item(15, 5).Value = 8091907.728333334
$ws.Cells.Item(15, 6).Value = 10114884.66041667
$ws.Cells.Item(15, 7).Value = 12137861.5925
$ws.Cells.Item(15, 8).Value = 14160838.52458333
$ws.Cells.Item(15, 9).Value = 16183815.45666667
$ws.Cells.Item(15, 10).Value = 18206792.38875
$ws.Cells.Item(15, 11).Value = 20229769.32083334
$ws.Cells.Item(15, 12).Value = 22252746.25291667
$ws.Cells.Item(15, 13).Value = 24275723.18500001
$ws.Cells.Item(15, 14).Value = 25788781.488750014
$ws.Cells.Item(15, 15).Value = 25788781.488750014
$ws.Cells.Item(15, 16).Value = 25788781.488750014
$ws.Cells.Item(15, 17).Value = 25788781.488750014
$ws.Cells.Item(15, 18).Value = 25788781.488750014
$ws.Cells.Item(15, 19).Value = 25788781.488750014
$ws.Cells.Item(15, 20).Value = 25788781.488750014
$ws.Cells.Item(15, 21).Value = 25788781.488750014
$ws.Cells.Item(15, 22).Value = 25788781.488750014
$ws.Cells.Item(15, 23).Value = 25788781.488750014
$ws.Cells.Item(15, 24).Value = 25788781.488750014
$ws.Cells.Item(15, 25).Value = 25788781.488750014
$ws.Cells.Item(15, 26).Value = 23765804.55666668
$ws.Cells.Item(15, 27).Value = 21742827.624583345
$ws.Cells.Item(15, 28).Value = 19719850.69250001
$ws.Cells.Item(15, 29).Value = 17696873.760416675
$ws.Cells.Item(15, 30).Value = 15673896.82833334
$ws.Cells.Item(15, 31).Value = 13650919.896250006
$ws.Cells.Item(15, 32).Value = 11627942.964166671
$ws.Cells.Item(15, 33).Value = 9604966.032083336
$ws.Cells.Item(15, 34).Value = 7581989.100000002
$ws.Cells.Item(15, 35).Value = 5559012.167916669
$ws.Cells.Item(15, 36).Value = 3536035.235833334
$ws.Cells.Item(15, 37).Value = 1513058.30375

# ===== Sheet: WIN =====
$ws = $wb.Worksheets.Item("WIN")

# New column header AK1 (date 20220218), matching style of existing header cells
$ws.Cells.Item(1, 37).NumberFormat = "@"
$ws.Cells.Item(1, 37).Font.Bold = $true
$ws.Cells.Item(1, 37).HorizontalAlignment = -4108
$ws.Cells.Item(1, 37).VerticalAlignment = -4160
$ws.Cells.Item(1, 37).Borders.LineStyle = 1
$ws.Cells.Item(1, 37).Value = "20220218"

# New row 14: period 20210903-20210910
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Font.Bold = $true
$ws.Cells.Item(14, 1).HorizontalAlignment = -4108
$ws.Cells.Item(14, 1).VerticalAlignment = -4160
$ws.Cells.Item(14, 1).Borders.LineStyle = 1
$ws.Cells.Item(14, 1).Value = "20210903-20210910"
$ws.Cells.Item(14, 14).Value = 3632703.442083333
$ws.Cells.Item(14, 15).Value = 3632703.442083333
$ws.Cells.Item(14, 16).Value = 3632703.442083333
$ws.Cells.Item(14, 17).Value = 3632703.442083333
$ws.Cells.Item(14, 18).Value = 3632703.442083333
$ws.Cells.Item(14, 19).Value = 3632703.442083333
$ws.Cells.Item(14, 20).Value = 3632703.442083333
$ws.Cells.Item(14, 21).Value = 3632703.442083333
$ws.Cells.Item(14, 22).Value = 3632703.442083333
$ws.Cells.Item(14, 23).Value = 3632703.442083333
$ws.Cells.Item(14, 24).Value = 3632703.442083333
$ws.Cells.Item(14, 25).Value = 3632703.442083333
$ws.Cells.Item(14, 26).Value = 3632703.442083333
$ws.Cells.Item(14, 27).Value = 3632703.442083333
$ws.Cells.Item(14, 28).Value = 3632703.442083333
$ws.Cells.Item(14, 29).Value = 3632703.442083333
$ws.Cells.Item(14, 30).Value = 3632703.442083333
$ws.Cells.Item(14, 31).Value = 3632703.442083333
$ws.Cells.Item(14, 32).Value = 3632703.442083333
$ws.Cells.Item(14, 33).Value = 3632703.442083333
$ws.Cells.Item(14, 34).Value = 3632703.442083333
$ws.Cells.Item(14, 35).Value = 3632703.442083333
$ws.Cells.Item(14, 36).Value = 3632703.442083333
$ws.Cells.Item(14, 37).Value = 3632703.442083333

# Row 15 (previously row 14): SUM row, shifted down, with new cumulative totals
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Font.Bold = $true
$ws.Cells.Item(15, 1).HorizontalAlignment = -4108
$ws.Cells.Item(15, 1).VerticalAlignment = -4160
$ws.Cells.Item(15, 1).Borders.LineStyle = 1
$ws.Cells.Item(15, 1).Value = "SUM"
$ws.Cells.Item(15, 2).Value = 5319315.755833333
$ws.Cells.Item(15, 3).Value = 10638631.51166667
$ws.Cells.Item(15, 4).Value = 15957947.2675
$ws.Cells.Item(15, 5).Value = 21277263.02333333
$ws.Cells.Item(15, 6).Value = 26596578.77916667
$ws.Cells.Item(15, 7).Value = 31915894.535
$ws.Cells.Item(15, 8).Value = 37235210.29083334
$ws.Cells.Item(15, 9).Value = 42554526.04666667
$ws.Cells.Item(15, 10).Value = 47873841.80250001
$ws.Cells.Item(15, 11).Value = 53193157.55833334
$ws.Cells.Item(15, 12).Value = 58512473.31416668
$ws.Cells.Item(15, 13).Value = 63831789.07000002
$ws.Cells.Item(15, 14).Value = 67464492.51208335
$ws.Cells.Item(15, 15).Value = 67464492.51208335
$ws.Cells.Item(15, 16).Value = 67464492.51208335
$ws.Cells.Item(15, 17).Value = 67464492.51208335
$ws.Cells.Item(15, 18).Value = 67464492.51208335
$ws.Cells.Item(15, 19).Value = 67464492.51208335
$ws.Cells.Item(15, 20).Value = 67464492.51208335
$ws.Cells.Item(15, 21).Value = 67464492.51208335
$ws.Cells.Item(15, 22).Value = 67464492.51208335
$ws.Cells.Item(15, 23).Value = 67464492.51208335
$ws.Cells.Item(15, 24).Value = 67464492.51208335
$ws.Cells.Item(15, 25).Value = 67464492.51208335
$ws.Cells.Item(15, 26).Value = 62145176.75625002
$ws.Cells.Item(15, 27).Value = 56825861.00041668
$ws.Cells.Item(15, 28).Value = 51506545.244583346
$ws.Cells.Item(15, 29).Value = 46187229.48875001
$ws.Cells.Item(15, 30).Value = 40867913.732916676
$ws.Cells.Item(15, 31).Value = 35548597.97708334
$ws.Cells.Item(15, 32).Value = 30229282.22125
$ws.Cells.Item(15, 33).Value = 24909966.465416666
$ws.Cells.Item(15, 34).Value = 19590650.709583335
$ws.Cells.Item(15, 35).Value = 14271334.95375
$ws.Cells.Item(15, 36).Value = 8952019.197916666
$ws.Cells.Item(15, 37).Value = 3632703.442083333

# ===== Sheet: NFT =====
$ws = $wb.Worksheets.Item("NFT")

# New column header AK1 (date 20220218), matching style of existing header cells
$ws.Cells.Item(1, 37).NumberFormat = "@"
$ws.Cells.Item(1, 37).Font.Bold = $true
$ws.Cells.Item(1, 37).HorizontalAlignment = -4108
$ws.Cells.Item(1, 37).VerticalAlignment = -4160
$ws.Cells.Item(1, 37).Borders.LineStyle = 1
$ws.Cells.Item(1, 37).Value = "20220218"

# New row 14: period 20210903-20210910
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Font.Bold = $true
$ws.Cells.Item(14, 1).HorizontalAlignment = -4108
$ws.Cells.Item(14, 1).VerticalAlignment = -4160
$ws.Cells.Item(14, 1).Borders.LineStyle = 1
$ws.Cells.Item(14, 1).Value = "20210903-20210910"
$ws.Cells.Item(14, 14).Value = 2257689781.774583
$ws.Cells.Item(14, 15).Value = 2257689781.774583
$ws.Cells.Item(14, 16).Value = 2257689781.774583
$ws.Cells.Item(14, 17).Value = 2257689781.774583
$ws.Cells.Item(14, 18).Value = 2257689781.774583
$ws.Cells.Item(14, 19).Value = 2257689781.774583
$ws.Cells.Item(14, 20).Value = 2257689781.774583
$ws.Cells.Item(14, 21).Value = 2257689781.774583
$ws.Cells.Item(14, 22).Value = 2257689781.774583
$ws.Cells.Item(14, 23).Value = 2257689781.774583
$ws.Cells.Item(14, 24).Value = 2257689781.774583
$ws.Cells.Item(14, 25).Value = 2257689781.774583
$ws.Cells.Item(14, 26).Value = 2257689781.774583
$ws.Cells.Item(14, 27).Value = 2257689781.774583
$ws.Cells.Item(14, 28).Value = 2257689781.774583
$ws.Cells.Item(14, 29).Value = 2257689781.774583
$ws.Cells.Item(14, 30).Value = 2257689781.774583
$ws.Cells.Item(14, 31).Value = 2257689781.774583
$ws.Cells.Item(14, 32).Value = 2257689781.774583
$ws.Cells.Item(14, 33).Value = 2257689781.774583
$ws.Cells.Item(14, 34).Value = 2257689781.774583
$ws.Cells.Item(14, 35).Value = 2257689781.774583
$ws.Cells.Item(14, 36).Value = 2257689781.774583
$ws.Cells.Item(14, 37).Value = 2257689781.774583

# Row 15 (previously row 14): SUM row, shifted down, with new cumulative totals
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Font.Bold = $true
$ws.Cells.Item(15, 1).HorizontalAlignment = -4108
$ws.Cells.Item(15, 1).VerticalAlignment = -4160
$ws.Cells.Item(15, 1).Borders.LineStyle = 1
$ws.Cells.Item(15, 1).Value = "SUM"
$ws.Cells.Item(15, 2).Value = 2833161790.824167
$ws.Cells.Item(15, 3).Value = 5666323581.648334
$ws.Cells.Item(15, 4).Value = 8499485372.4725
$ws.Cells.Item(15, 5).Value = 11332647163.29667
$ws.Cells.Item(15, 6).Value = 14165808954.12083
$ws.Cells.Item(15, 7).Value = 16998970744.945
$ws.Cells.Item(15, 8).Value = 19832132535.76917
$ws.Cells.Item(15, 9).Value = 22665294326.59333
$ws.Cells.Item(15, 10).Value = 25498456117.4175
$ws.Cells.Item(15, 11).Value = 28331617908.24166
$ws.Cells.Item(15, 12).Value = 31164779699.06583
$ws.Cells.Item(15, 13).Value = 33997941489.89
$ws.Cells.Item(15, 14).Value = 36255631271.66458
$ws.Cells.Item(15, 15).Value = 36255631271.66458
$ws.Cells.Item(15, 16).Value = 36255631271.66458
$ws.Cells.Item(15, 17).Value = 36255631271.66458
$ws.Cells.Item(15, 18).Value = 36255631271.66458
$ws.Cells.Item(15, 19).Value = 36255631271.66458
$ws.Cells.Item(15, 20).Value = 36255631271.66458
$ws.Cells.Item(15, 21).Value = 36255631271.66458
$ws.Cells.Item(15, 22).Value = 36255631271.66458
$ws.Cells.Item(15, 23).Value = 36255631271.66458
$ws.Cells.Item(15, 24).Value = 36255631271.66458
$ws.Cells.Item(15, 25).Value = 36255631271.66458
$ws.Cells.Item(15, 26).Value = 33422469480.840412
$ws.Cells.Item(15, 27).Value = 30589307690.016247
$ws.Cells.Item(15, 28).Value = 27756145899.19208
$ws.Cells.Item(15, 29).Value = 24922984108.367916
$ws.Cells.Item(15, 30).Value = 22089822317.54375
$ws.Cells.Item(15, 31).Value = 19256660526.719585
$ws.Cells.Item(15, 32).Value = 16423498735.895416
$ws.Cells.Item(15, 33).Value = 13590336945.07125
$ws.Cells.Item(15, 34).Value = 10757175154.247082
$ws.Cells.Item(15, 35).Value = 7924013363.422916
$ws.Cells.Item(15, 36).Value = 5090851572.598749
$ws.Cells.Item(15, 37).Value = 2257689781.774583
